$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1325.6885
$ws.Range("I17").Value = 499
$ws.Range("J17").Value = 1368.4482
$ws.Range("K17").Value = 1497
$ws.Range("L17").Value = 4105.3446
$ws.Range("M17").Value = -1329
$ws.Range("N17").Value = -4441.3446

$ws.Range("H135").Value = 331.35715
$ws.Range("I135").Value = 331.35715
$ws.Range("K135").Value = 2982.21435
$ws.Range("M135").Value = -447.2143499999997


# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 1530.4
$ws.Range("I4").Value = 1900.5714
$ws.Range("J4").Value = 666.6667
$ws.Range("K4").Value = 1900.5714
$ws.Range("L4").Value = 666.6667
$ws.Range("M4").Value = -1784.5714
$ws.Range("N4").Value = -898.6667

$ws.Range("H5").Value = 195.5
$ws.Range("I5").Value = 150.71428
$ws.Range("J5").Value = 300
$ws.Range("K5").Value = 150.71428
$ws.Range("L5").Value = 300
$ws.Range("M5").Value = -38.71428
$ws.Range("N5").Value = -524

$ws.Range("H14").Value = 722222200
$ws.Range("I14").Value = 444444450
$ws.Range("K14").Value = 444444450
$ws.Range("M14").Value = -444444275

$ws.Range("H21").Value = 2345
$ws.Range("I21").Value = 2345
$ws.Range("K21").Value = 2345
$ws.Range("M21").Value = -1971

$ws.Range("H22").Value = 6623.2
$ws.Range("I22").Value = 1608
$ws.Range("J22").Value = 9966.666999999999
$ws.Range("K22").Value = 1608
$ws.Range("L22").Value = 9966.666999999999
$ws.Range("M22").Value = -1309
$ws.Range("N22").Value = -10564.667

$ws.Range("H26").Value = 0
$ws.Range("I26").Value = 0
$ws.Range("J26").Value = 0
$ws.Range("K26").Value = 0
$ws.Range("L26").Value = 0
$ws.Range("M26").ClearContents()
$ws.Range("N26").ClearContents()

$ws.Range("H27").Value = 6000
$ws.Range("J27").Value = 6000
$ws.Range("L27").Value = 6000
$ws.Range("N27").Value = -6368

$ws.Range("H33").Value = 0
$ws.Range("I33").Value = 0
$ws.Range("K33").Value = 0
$ws.Range("M33").ClearContents()

$ws.Range("H38").Value = 0
$ws.Range("J38").Value = 0
$ws.Range("L38").Value = 0
$ws.Range("N38").ClearContents()

$ws.Range("H39").Value = 0
$ws.Range("I39").Value = 0
$ws.Range("K39").Value = 0
$ws.Range("M39").ClearContents()

$ws.Range("H40").Value = 0
$ws.Range("I40").Value = 0
$ws.Range("K40").Value = 0
$ws.Range("M40").ClearContents()

$ws.Range("H51").Value = 0
$ws.Range("J51").Value = 0
$ws.Range("L51").Value = 0
$ws.Range("N51").ClearContents()

$ws.Range("H58").Value = 0
$ws.Range("J58").Value = 0
$ws.Range("L58").Value = 0
$ws.Range("N58").ClearContents()

$ws.Range("H61").Value = 401870.44
$ws.Range("I61").Value = 2041.4706
$ws.Range("K61").Value = 2041.4706
$ws.Range("M61").Value = -1829.4706

$ws.Range("H136").Value = 401870.44
$ws.Range("I136").Value = 2041.4706
$ws.Range("K136").Value = 6124.4118
$ws.Range("M136").Value = -3574.4118


# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 195.5
$ws.Range("I4").Value = 150.71428
$ws.Range("J4").Value = 300
$ws.Range("K4").Value = 150.71428
$ws.Range("L4").Value = 300
$ws.Range("M4").Value = -35.71428
$ws.Range("N4").Value = -530

$ws.Range("H10").Value = 0
$ws.Range("J10").Value = 0
$ws.Range("L10").Value = 0
$ws.Range("N10").ClearContents()

$ws.Range("H22").Value = 155.88889
$ws.Range("I22").Value = 138.5
$ws.Range("J22").Value = 295
$ws.Range("K22").Value = 138.5
$ws.Range("L22").Value = 295
$ws.Range("M22").Value = 34.5
$ws.Range("N22").Value = -641

$ws.Range("H57").Value = 0
$ws.Range("J57").Value = 0
$ws.Range("L57").Value = 0
$ws.Range("N57").ClearContents()

$ws.Range("H58").Value = 48125
$ws.Range("I58").Value = 0
$ws.Range("J58").Value = 48125
$ws.Range("K58").Value = 0
$ws.Range("L58").Value = 48125
$ws.Range("M58").ClearContents()
$ws.Range("N58").Value = -48713

$ws.Range("H59").Value = 50000
$ws.Range("J59").Value = 50000
$ws.Range("L59").Value = 50000
$ws.Range("N59").Value = -51694

$ws.Range("H60").Value = 40000
$ws.Range("I60").Value = 0
$ws.Range("J60").Value = 40000
$ws.Range("K60").Value = 0
$ws.Range("L60").Value = 40000
$ws.Range("M60").ClearContents()
$ws.Range("N60").Value = -41198

$ws.Range("H136").Value = 0
$ws.Range("J136").Value = 0
$ws.Range("L136").Value = 0
$ws.Range("N136").ClearContents()


# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 5658.7144
$ws.Range("I16").Value = 3322.2
$ws.Range("J16").Value = 11500
$ws.Range("K16").Value = 3322.2
$ws.Range("L16").Value = 11500
$ws.Range("M16").Value = -3035.2
$ws.Range("N16").Value = -12074

$ws.Range("H52").Value = 19225.715
$ws.Range("J52").Value = 19225.715
$ws.Range("L52").Value = 19225.715
$ws.Range("N52").Value = -19813.715

$ws.Range("H58").Value = 1638.4762
$ws.Range("I58").Value = 1309
$ws.Range("J58").Value = 1803.2142
$ws.Range("K58").Value = 1309
$ws.Range("L58").Value = 1803.2142
$ws.Range("M58").Value = -1106
$ws.Range("N58").Value = -2209.2142

$ws.Range("H105").Value = 20000
$ws.Range("I105").Value = 0
$ws.Range("J105").Value = 20000
$ws.Range("K105").Value = 0
$ws.Range("L105").Value = 20000
$ws.Range("M105").ClearContents()
$ws.Range("N105").Value = -23494

$ws.Range("H113").Value = 5658.7144
$ws.Range("I113").Value = 3322.2
$ws.Range("J113").Value = 11500
$ws.Range("K113").Value = 3322.2
$ws.Range("L113").Value = 11500
$ws.Range("M113").Value = -1152.2
$ws.Range("N113").Value = -15840

$ws.Range("H136").Value = 1638.4762
$ws.Range("I136").Value = 1309
$ws.Range("J136").Value = 1803.2142
$ws.Range("K136").Value = 3927
$ws.Range("L136").Value = 5409.642599999999
$ws.Range("M136").Value = -1377
$ws.Range("N136").Value = -10509.6426


# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 843.7273
$ws.Range("I2").Value = 1019.8889
$ws.Range("J2").Value = 51
$ws.Range("K2").Value = 6119.3334
$ws.Range("L2").Value = 306
$ws.Range("M2").Value = -6006.3334
$ws.Range("N2").Value = -532

$ws.Range("H9").Value = 1501500.5
$ws.Range("J9").Value = 3000
$ws.Range("L9").Value = 9000
$ws.Range("N9").Value = -9448

$ws.Range("H15").Value = 345.14285
$ws.Range("I15").Value = 211.2
$ws.Range("J15").Value = 680
$ws.Range("K15").Value = 633.5999999999999
$ws.Range("L15").Value = 2040
$ws.Range("M15").Value = -493.5999999999999
$ws.Range("N15").Value = -2320

$ws.Range("H20").Value = 4250.4
$ws.Range("I20").Value = 1500
$ws.Range("J20").Value = 4938
$ws.Range("K20").Value = 4500
$ws.Range("L20").Value = 14814
$ws.Range("M20").Value = -4273
$ws.Range("N20").Value = -15268

$ws.Range("H21").Value = 725
$ws.Range("I21").Value = 450
$ws.Range("J21").Value = 1000
$ws.Range("K21").Value = 1350
$ws.Range("L21").Value = 3000
$ws.Range("M21").Value = -1177
$ws.Range("N21").Value = -3346

$ws.Range("H22").Value = 55890890
$ws.Range("J22").Value = 66669068
$ws.Range("L22").Value = 200007204
$ws.Range("N22").Value = -200007542

$ws.Range("H26").Value = 1081.5
$ws.Range("I26").Value = 96.875
$ws.Range("J26").Value = 5020
$ws.Range("K26").Value = 290.625
$ws.Range("L26").Value = 15060
$ws.Range("M26").Value = -2.625
$ws.Range("N26").Value = -15636

$ws.Range("H27").Value = 55890890
$ws.Range("J27").Value = 66669068
$ws.Range("L27").Value = 200007204
$ws.Range("N27").Value = -200007408


# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H137").Value = 39991.43
$ws.Range("I137").Value = 5000
$ws.Range("J137").Value = 45823.332
$ws.Range("K137").Value = 5000
$ws.Range("L137").Value = 45823.332
$ws.Range("M137").Value = 100
$ws.Range("N137").Value = -56023.332


# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 25732.25
$ws.Range("I132").Value = 5542.7144
$ws.Range("K132").Value = 16628.1432
$ws.Range("M132").Value = -14098.1432

$ws.Range("H136").Value = 4222.59
$ws.Range("I136").Value = 1460.7307
$ws.Range("J136").Value = 9746.308000000001
$ws.Range("K136").Value = 4382.1921
$ws.Range("L136").Value = 29238.924
$ws.Range("M136").Value = -1832.1921
$ws.Range("N136").Value = -34338.924


# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1907.9166
$ws.Range("I132").Value = 1980.7
$ws.Range("K132").Value = 5942.1
$ws.Range("M132").Value = -3412.1

$ws.Range("H136").Value = 3152.4285
$ws.Range("I136").Value = 3435.7026
$ws.Range("J136").Value = 2600.7896
$ws.Range("K136").Value = 10307.1078
$ws.Range("L136").Value = 7802.3688
$ws.Range("M136").Value = -7757.1078
$ws.Range("N136").Value = -12902.3688

